$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The single "AmazonLogin" test case row is being split into two rows:
# one for an incorrect login attempt and one for a correct login attempt.
# Duplicate row 2 (copy, so formatting/styles carry over) and insert the
# copy right below it; this pushes the old rows 3 and 4
# (AmazonSearchItem / AmazonProceedToCheckOut) down to rows 4 and 5.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# Row 2: the "incorrect login" test case.
$ws.Range("A2").Value = "AmazonIncorrectLogin"
$ws.Range("B2").Value = 7020214690
$ws.Range("C2").Value = "test12"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"

# Row 3: the "correct login" test case (what used to be plain "AmazonLogin").
$ws.Range("A3").Value = "AmazonCorrectLogin"
$ws.Range("B3").Value = 7020214690
$ws.Range("C3").Value = "test123"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"

# Rows 4 and 5 (previously rows 3 and 4) already hold the correct
# AmazonSearchItem / AmazonProceedToCheckOut contents, untouched.

# Update the selection to match the new active cell after editing.
$ws.Range("D6").Select() | Out-Null
